# Add "MicroVG" and "GPU" to the GUI overview schema
#
# - "Drawing"  -> "MicroVG"        (Rounded Rectangle 36)
# - "Touch"    -> "Input devices"  (Rounded Rectangle 38)
# - "Buttons"  -> "GPU"            (Rounded Rectangle 39)
# - date placeholders "janvier 23" -> "octobre 23" in the handout / notes masters

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Slide shape text updates -------------------------------------------------

$drawingShape = $s.Shapes.Item("Rounded Rectangle 36")
$drawingShape.TextFrame.TextRange.Text = "MicroVG"

$touchShape = $s.Shapes.Item("Rounded Rectangle 38")
$touchShape.TextFrame.TextRange.Text = "Input devices"

$buttonsShape = $s.Shapes.Item("Rounded Rectangle 39")
$buttonsShape.TextFrame.TextRange.Text = "GPU"

# --- Footer / date placeholders -----------------------------------------------
# The handout master and notes master both show a fixed-date field reading
# "janvier 23" - update them to "octobre 23" via the HeadersFooters object.

$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.UseFormat = $false
$handoutDate.Value = "octobre 23"

$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.UseFormat = $false
$notesDate.Value = "octobre 23"
